$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "90.573.54"
Set-TextCell "E2" "  -0.02%  "
Set-TextCell "D3" "3.100.80"
Set-TextCell "E3" "  -0.92%  "
Set-TextCell "E4" "  +0.11%  "
Set-TextCell "D5" "241.65"
Set-TextCell "E5" "  +2.84%  "
Set-TextCell "D6" "625.63"
Set-TextCell "E6" "  -0.06%  "
Set-TextCell "E7" "  +8.24%  "
Set-TextCell "D8" "0.373"
Set-TextCell "E8" "  +4.52%  "
Set-TextCell "E9" "  +0.12%  "
Set-TextCell "D10" "0.738"
Set-TextCell "E10" "  +2.45%  "
Set-TextCell "E11" "  -11.96%  "
Set-TextCell "E12" "  +3.71%  "
Set-TextCell "D13" "0.0000249"
Set-TextCell "E13" "  +2.39%  "
Set-TextCell "D14" "35.26"
Set-TextCell "E14" "  -3.02%  "
Set-TextCell "D15" "5.49"
Set-TextCell "E15" "  -2.02%  "
Set-TextCell "D16" "90.472.30"
Set-TextCell "E16" "  +0.33%  "
Set-TextCell "D17" "3.677.05"
Set-TextCell "E17" "  +0.10%  "
Set-TextCell "D18" "3.150.01"
Set-TextCell "E18" "  +2.67%  "
Set-TextCell "D19" "3.85"
Set-TextCell "E19" "  +4.43%  "
Set-TextCell "D20" "14.24"
Set-TextCell "E20" "  -1.05%  "
Set-TextCell "D21" "0.0000210"
Set-TextCell "E21" "  -1.25%  "
Set-TextCell "D22" "5.77"
Set-TextCell "E22" "  +6.94%  "
Set-TextCell "D23" "445.95"
Set-TextCell "E23" "  -0.88%  "
Set-TextCell "D24" "9.09"
Set-TextCell "E24" "  +0.53%  "
Set-TextCell "D25" "5.90"
Set-TextCell "E25" "  +0.47%  "
Set-TextCell "D26" "92.92"
Set-TextCell "E26" "  +2.34%  "
Set-TextCell "D27" "12.04"
Set-TextCell "E27" "  -1.61%  "
Set-TextCell "E28" "  +0.24%  "
Set-TextCell "E29" "  +0.13%  "
Set-TextCell "D30" "0.175"
Set-TextCell "E30" "  +10.07%  "
Set-TextCell "D31" "9.26"
Set-TextCell "E31" "  -0.05%  "
Set-TextCell "D32" "0.217"
Set-TextCell "E32" "  +11.13%  "
Set-TextCell "D33" "0.998"
Set-TextCell "E33" "  +6.35%  "
Set-TextCell "D34" "0.108"
Set-TextCell "E34" "  +27.71%  "
Set-TextCell "D35" "4.40"
Set-TextCell "E35" "  +37.83%  "
Set-TextCell "D36" "26.59"
Set-TextCell "E36" "  -4.26%  "
Set-TextCell "B37" "Kaspa"
Set-TextCell "C37" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D37" "0.157"
Set-TextCell "E37" "  +3.78%  "
Set-TextCell "B38" "RenderToken"
Set-TextCell "C38" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextCell "D38" "7.54"
Set-TextCell "E38" "  +7.66%  "
Set-TextCell "B39" "PancakeSwap"
Set-TextCell "C39" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D39" "1.92"
Set-TextCell "E39" "  +0.11%  "
Set-TextCell "B40" "Bittensor"
Set-TextCell "C40" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D40" "494.50"
Set-TextCell "E40" "  -2.44%  "
Set-TextCell "D41" "3.60"
Set-TextCell "E41" "  -1.93%  "
Set-TextCell "D42" "1.29"
Set-TextCell "E42" "  -1.78%  "
Set-TextCell "E43" "  -1.98%  "
Set-TextCell "D44" "22.12"
Set-TextCell "E44" "  -0.32%  "
Set-TextCell "E45" "  +0.02%  "
Set-TextCell "D46" "159.34"
Set-TextCell "E46" "  +6.96%  "
Set-TextCell "D47" "1.90"
Set-TextCell "E47" "  -3.14%  "
Set-TextCell "D48" "0.686"
Set-TextCell "E48" "  -2.11%  "
Set-TextCell "D49" "4.56"
Set-TextCell "E49" "  -0.24%  "
Set-TextCell "D50" "45.04"
Set-TextCell "E50" "  +0.97%  "
Set-TextCell "D51" "1.34"
Set-TextCell "E51" "  -1.30%  "
